# Add more responses to survey
# Applies new Microsoft Forms survey responses to the "English" and
# "Chinese" sheets of the workbook, plus the formatting side-effects
# that Excel applied when the extra rows were appended (a handful of
# pre-existing cells lost their "vertical center" style in favour of
# plain "wrap text" once the sheet was re-saved, and several rows grew
# taller to fit newly-entered, longer answers).

$wb = $excel.ActiveWorkbook
$wsEnglish = $wb.Worksheets.Item("English")
$wsChinese = $wb.Worksheets.Item("Chinese")

# ---------------------------------------------------------------------
# 1. English sheet: style + row-height touch-ups on pre-existing rows
# ---------------------------------------------------------------------
# U4/U5/U10/V15/V37/U51/V61/V65 move from style index 3 (vertical-center)
# to style index 2 (wrap-text) -- copy the format from a neighbouring
# cell that already carries style 2 so no new style entries are created.
$wsEnglish.Range("T4").Copy()
$wsEnglish.Range("U4").PasteSpecial(-4122)

$wsEnglish.Range("T5").Copy()
$wsEnglish.Range("U5").PasteSpecial(-4122)

$wsEnglish.Range("T10").Copy()
$wsEnglish.Range("U10").PasteSpecial(-4122)

$wsEnglish.Range("U15").Copy()
$wsEnglish.Range("V15").PasteSpecial(-4122)

$wsEnglish.Range("U37").Copy()
$wsEnglish.Range("V37").PasteSpecial(-4122)

$wsEnglish.Range("T51").Copy()
$wsEnglish.Range("U51").PasteSpecial(-4122)

$wsEnglish.Range("U61").Copy()
$wsEnglish.Range("V61").PasteSpecial(-4122)

$wsEnglish.Range("U65").Copy()
$wsEnglish.Range("V65").PasteSpecial(-4122)

# Row heights grew to fit the (re-wrapped) existing text.
$wsEnglish.Rows.Item(4).RowHeight = 123.9
$wsEnglish.Rows.Item(5).RowHeight = 123.9
$wsEnglish.Rows.Item(51).RowHeight = 197.7
$wsEnglish.Rows.Item(65).RowHeight = 123.9

# ---------------------------------------------------------------------
# 2. English sheet: four new survey responses (rows 76-79)
# ---------------------------------------------------------------------
# Row 4 now has the canonical formatting pattern for a response row
# (A/G:O right-aligned numeric style, everything else wrap-text) -- use
# it as the paste-format template for the freshly appended rows.
$wsEnglish.Range("A4:W4").Copy()
$wsEnglish.Range("A76:W76").PasteSpecial(-4122)
$wsEnglish.Range("A4:W4").Copy()
$wsEnglish.Range("A77:W77").PasteSpecial(-4122)
$wsEnglish.Range("A4:W4").Copy()
$wsEnglish.Range("A78:W78").PasteSpecial(-4122)
$wsEnglish.Range("A4:W4").Copy()
$wsEnglish.Range("A79:W79").PasteSpecial(-4122)

# Row 76
$wsEnglish.Cells.Item(76, 1).Value = "23/04/2024 20:16:50"
$wsEnglish.Cells.Item(76, 2).Value = "United States"
$wsEnglish.Cells.Item(76, 3).Value = "No"
$wsEnglish.Cells.Item(76, 4).Value = "30-44"
$wsEnglish.Cells.Item(76, 5).Value = "Male"
$wsEnglish.Cells.Item(76, 6).Value = "Desktop / Laptop"
$wsEnglish.Cells.Item(76, 7).Value = 5
$wsEnglish.Cells.Item(76, 8).Value = 5
$wsEnglish.Cells.Item(76, 9).Value = 5
$wsEnglish.Cells.Item(76, 10).Value = 5
$wsEnglish.Cells.Item(76, 11).Value = 3
$wsEnglish.Cells.Item(76, 12).Value = 5
$wsEnglish.Cells.Item(76, 13).Value = 3
$wsEnglish.Cells.Item(76, 14).Value = 5
$wsEnglish.Cells.Item(76, 15).Value = 3
$wsEnglish.Cells.Item(76, 16).Value = "Wu, Jian-Hong (吳儉鴻）"
$wsEnglish.Cells.Item(76, 17).Value = 'I feel that I can relate to Mr. WU Jian-Hong as I am also an elementary school teacher. I love that he has taught for so many years and sounds like he is leaving behind a wonderful legacy!'
$wsEnglish.Cells.Item(76, 18).Value = "Shilin Elementary School (士林國小)"
$wsEnglish.Cells.Item(76, 19).Value = "Being an educator, I enjoyed learning about the school."
$wsEnglish.Cells.Item(76, 20).Value = "How easy it was to navigate it."
$wsEnglish.Cells.Item(76, 21).Value = "I did not find anything I disliked about it."
$wsEnglish.Cells.Item(76, 22).Value = "I thought it was great."
$wsEnglish.Cells.Item(76, 23).Value = "Well done!"
$wsEnglish.Rows.Item(76).RowHeight = 271.5

# Row 77
$wsEnglish.Cells.Item(77, 1).Value = "23/04/2024 23:18:43"
$wsEnglish.Cells.Item(77, 2).Value = "United States"
$wsEnglish.Cells.Item(77, 3).Value = "No"
$wsEnglish.Cells.Item(77, 4).Value = "18-29"
$wsEnglish.Cells.Item(77, 5).Value = "Female"
$wsEnglish.Cells.Item(77, 6).Value = "Smartphone"
$wsEnglish.Cells.Item(77, 7).Value = 5
$wsEnglish.Cells.Item(77, 8).Value = 5
$wsEnglish.Cells.Item(77, 9).Value = 5
$wsEnglish.Cells.Item(77, 10).Value = 5
$wsEnglish.Cells.Item(77, 11).Value = 5
$wsEnglish.Cells.Item(77, 12).Value = 5
$wsEnglish.Cells.Item(77, 13).Value = 3
$wsEnglish.Cells.Item(77, 14).Value = 5
$wsEnglish.Cells.Item(77, 15).Value = 4
$wsEnglish.Cells.Item(77, 16).Value = "Lily (莉莉)"
$wsEnglish.Cells.Item(77, 17).Value = "It was interesting reading about the traditions behind her current work"
$wsEnglish.Cells.Item(77, 18).Value = "Zhishanyan Huiji Temple (芝山巖惠濟宮)"
$wsEnglish.Cells.Item(77, 20).Value = "Very aesthetically pleasing"
$wsEnglish.Cells.Item(77, 21).Value = "There wasn’t anything"
$wsEnglish.Rows.Item(77).RowHeight = 111.6

# Row 78
$wsEnglish.Cells.Item(78, 1).Value = "24/04/2024 01:36:05"
$wsEnglish.Cells.Item(78, 2).Value = "United States"
$wsEnglish.Cells.Item(78, 3).Value = "No"
$wsEnglish.Cells.Item(78, 4).Value = "30-44"
$wsEnglish.Cells.Item(78, 5).Value = "Female"
$wsEnglish.Cells.Item(78, 6).Value = "Desktop / Laptop"
$wsEnglish.Cells.Item(78, 7).Value = 5
$wsEnglish.Cells.Item(78, 8).Value = 5
$wsEnglish.Cells.Item(78, 9).Value = 5
$wsEnglish.Cells.Item(78, 10).Value = 5
$wsEnglish.Cells.Item(78, 11).Value = 5
$wsEnglish.Cells.Item(78, 12).Value = 5
$wsEnglish.Cells.Item(78, 13).Value = 5
$wsEnglish.Cells.Item(78, 14).Value = 5
$wsEnglish.Cells.Item(78, 15).Value = 5
$wsEnglish.Cells.Item(78, 16).Value = "Lily (莉莉)"
$wsEnglish.Cells.Item(78, 18).Value = "Shilin Elementary School (士林國小)"
$wsEnglish.Cells.Item(78, 20).Value = "I loved the opening page"
$wsEnglish.Rows.Item(78).RowHeight = 62.4

# Row 79
$wsEnglish.Cells.Item(79, 1).Value = "24/04/2024 23:23:44"
$wsEnglish.Cells.Item(79, 2).Value = "United States"
$wsEnglish.Cells.Item(79, 3).Value = "No"
$wsEnglish.Cells.Item(79, 4).Value = "18-29"
$wsEnglish.Cells.Item(79, 5).Value = "Female"
$wsEnglish.Cells.Item(79, 6).Value = "Smartphone"
$wsEnglish.Cells.Item(79, 7).Value = 5
$wsEnglish.Cells.Item(79, 8).Value = 5
$wsEnglish.Cells.Item(79, 9).Value = 5
$wsEnglish.Cells.Item(79, 10).Value = 5
$wsEnglish.Cells.Item(79, 11).Value = 4
$wsEnglish.Cells.Item(79, 12).Value = 3
$wsEnglish.Cells.Item(79, 13).Value = 5
$wsEnglish.Cells.Item(79, 14).Value = 4
$wsEnglish.Cells.Item(79, 15).Value = 5
$wsEnglish.Cells.Item(79, 16).Value = "Lily (莉莉)"
$wsEnglish.Cells.Item(79, 17).Value = 'My eye was caught immediately when her description said she was a descendant of one of the oldest families of Shilin, and the article did not disappoint. Although I have never visited, her family’s pharmacy is such an important place marker when considering modernity and the switch to Western medicine that her and her family as well as the country undertook during a period of modernization. It was such an important story to share, and I am so happy that I got to learn about her via this website!'
$wsEnglish.Cells.Item(79, 18).Value = "Zhishanyan Huiji Temple (芝山巖惠濟宮)"
$wsEnglish.Cells.Item(79, 19).Value = 'Because so much of this website focuses on the timeline and the circumstances under which the community and country made the transition into more modern or western concepts, including the temple as a marker of tradition even in a modern world is pivotal. Even though so much has changed in the centuries that this district has existed for, the permanence of the temple as both a material concept as well as traditional one shows the importance of sacred traditions even in a technologically changing world.'
$wsEnglish.Cells.Item(79, 20).Value = 'It’s very fluid. It was very easy for me to grasp the concept of the website (it almost flows like an essay in which you present traditional histories/practices in conjunction with modern life), and I really really enjoyed the personal stories that were shared in the “People” section.'
$wsEnglish.Cells.Item(79, 21).Value = "No notes!"
$wsEnglish.Cells.Item(79, 22).Value = "Maybe more stories! Stories from younger people perhaps? Not really sure!"
$wsEnglish.Cells.Item(79, 23).Value = 'Amazing job! Will be sharing this with my family and friends because it is great work! I can’t wait to visit one day!!'
$wsEnglish.Rows.Item(79).RowHeight = 409.6

# W79 keeps the "vertical center" style (index 3), matching the last
# pre-existing response row rather than the new wrap-text template.
$wsEnglish.Range("W75").Copy()
$wsEnglish.Cells.Item(79, 23).PasteSpecial(-4122)
$wsEnglish.Cells.Item(79, 23).Value = 'Amazing job! Will be sharing this with my family and friends because it is great work! I can’t wait to visit one day!!'

# ---------------------------------------------------------------------
# 3. Chinese sheet: two new survey responses (rows 53-54)
# ---------------------------------------------------------------------
$wsChinese.Range("A52:W52").Copy()
$wsChinese.Range("A53:W53").PasteSpecial(-4122)
$wsChinese.Range("A52:W52").Copy()
$wsChinese.Range("A54:W54").PasteSpecial(-4122)

# Row 53
$wsChinese.Cells.Item(53, 1).Value = "23/04/2024 15:54:41"
$wsChinese.Cells.Item(53, 2).Value = "臺灣"
$wsChinese.Cells.Item(53, 3).Value = "不"
$wsChinese.Cells.Item(53, 4).Value = "18-29"
$wsChinese.Cells.Item(53, 5).Value = "女性"
$wsChinese.Cells.Item(53, 6).Value = "手機"
$wsChinese.Cells.Item(53, 7).Value = 4
$wsChinese.Cells.Item(53, 8).Value = 4
$wsChinese.Cells.Item(53, 9).Value = 5
$wsChinese.Cells.Item(53, 10).Value = 5
$wsChinese.Cells.Item(53, 11).Value = 5
$wsChinese.Cells.Item(53, 12).Value = 4
$wsChinese.Cells.Item(53, 13).Value = 3
$wsChinese.Cells.Item(53, 14).Value = 4
$wsChinese.Cells.Item(53, 15).Value = 4
$wsChinese.Cells.Item(53, 16).Value = "Lily (莉莉)"
$wsChinese.Cells.Item(53, 17).Value = 'I’m very impressed with the story that she was experienced in person.'
$wsChinese.Cells.Item(53, 18).Value = "Shilin Architecture (士林建築)"
$wsChinese.Cells.Item(53, 19).Value = 'I love the words that you use to describe the dry river. It’s very happy to see the trace of the old river still remaining there.'
$wsChinese.Cells.Item(53, 20).Value = 'The website’s completion, the stories, and the way you introduce the story and history.'
$wsChinese.Cells.Item(53, 21).Value = "I think there are still some things that need to be completed."
$wsChinese.Cells.Item(53, 23).Value = 'You guys did a great job! I’m very surprised that how well the website is performing.'
$wsChinese.Rows.Item(53).RowHeight = 185.4

# U53 and W53 use the "vertical center" style (index 3) -- copy it from
# an existing style-3 cell elsewhere on the sheet (the style, not the
# column, is what PasteSpecial(formats) carries over).
$wsChinese.Range("R12").Copy()
$wsChinese.Cells.Item(53, 21).PasteSpecial(-4122)
$wsChinese.Cells.Item(53, 21).Value = "I think there are still some things that need to be completed."
$wsChinese.Range("R12").Copy()
$wsChinese.Cells.Item(53, 23).PasteSpecial(-4122)
$wsChinese.Cells.Item(53, 23).Value = 'You guys did a great job! I’m very surprised that how well the website is performing.'

# Row 54
$wsChinese.Cells.Item(54, 1).Value = "24/04/2024 23:06:51"
$wsChinese.Cells.Item(54, 2).Value = "臺灣"
$wsChinese.Cells.Item(54, 3).Value = "不"
$wsChinese.Cells.Item(54, 4).Value = "18-29"
$wsChinese.Cells.Item(54, 5).Value = "女性"
$wsChinese.Cells.Item(54, 6).Value = "手機"
$wsChinese.Cells.Item(54, 7).Value = 4
$wsChinese.Cells.Item(54, 8).Value = 5
$wsChinese.Cells.Item(54, 9).Value = 5
$wsChinese.Cells.Item(54, 10).Value = 4
$wsChinese.Cells.Item(54, 11).Value = 4
$wsChinese.Cells.Item(54, 12).Value = 4
$wsChinese.Cells.Item(54, 13).Value = 4
$wsChinese.Cells.Item(54, 14).Value = 4
$wsChinese.Cells.Item(54, 15).Value = 3
$wsChinese.Cells.Item(54, 16).Value = "Lily (莉莉)"
$wsChinese.Cells.Item(54, 18).Value = "Shilin Elementary School (士林國小)"
$wsChinese.Cells.Item(54, 20).Value = "圖片很漂亮，尤其是那張士林夜景圖"
$wsChinese.Cells.Item(54, 22).Value = "在人名旁邊可以加上一句話介紹，讓瀏覽網頁的人可以先簡單知道這個人的身分，在選擇的時候也更有方向性！"
$wsChinese.Rows.Item(54).RowHeight = 50.1

# R54 and V54 use the "vertical center" style (index 3).
$wsChinese.Range("R12").Copy()
$wsChinese.Cells.Item(54, 18).PasteSpecial(-4122)
$wsChinese.Cells.Item(54, 18).Value = "Shilin Elementary School (士林國小)"
$wsChinese.Range("R12").Copy()
$wsChinese.Cells.Item(54, 22).PasteSpecial(-4122)
$wsChinese.Cells.Item(54, 22).Value = "在人名旁邊可以加上一句話介紹，讓瀏覽網頁的人可以先簡單知道這個人的身分，在選擇的時候也更有方向性！"

# ---------------------------------------------------------------------
# 4. Selection / view state
# ---------------------------------------------------------------------
# Chinese keeps a "select everything" view (now covering the two new
# rows); English ends up scrolled down with a single active cell. Select
# English last so it remains the active / tabSelected sheet.
$wsChinese.Range("A1:W54").Select()
$wsEnglish.Activate()
$wsEnglish.Range("A3").Select()
